# Add new customer row: phone 79174421, blank birthday, 0 total_points
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage the phone number as text in a scratch cell, then copy/paste-values it
# into A7 so it lands as a genuine text cell (no quote-prefix/autodetect
# number conversion, matching how this phone column already stores values).
$scratch = $ws.Range("E1")
$scratch.NumberFormat = "@"
$scratch.Value = "79174421"
$scratch.Copy()
$ws.Range("A7").PasteSpecial(-4163)
$scratch.Clear()

# New customer has no birthday on file (matches the blank cells in B4:B6)
$ws.Range("B7").Value = ""

# Points reset to 0
$ws.Range("C7").Value = 0
